# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial 45171 (2023-09-02) to serial 45172 (2023-09-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 264; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45171) {
        $cell.Value2 = 45172
    }
}
